# Reverse the order of comma-separated entries ("Recorded By") in column G
# for every data row in the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $text.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $reversed = $trimmed[($trimmed.Count - 1)..0]
    $newText = [string]::Join(", ", $reversed)

    if ($newText -ne $text) {
        $cell.Value2 = $newText
    }
}
